$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "model_2_8_0"
$ws.Cells.Item(2, 2).Value = 0.9784118741222514
$ws.Cells.Item(2, 3).Value = 0.9981132030219909
$ws.Cells.Item(2, 4).Value = 0.9966477698062028
$ws.Cells.Item(2, 5).Value = 0.9976767746371612
$ws.Cells.Item(2, 6).Value = 2.603534408302417
$ws.Cells.Item(2, 7).Value = 0.2876527593526298
$ws.Cells.Item(2, 8).Value = 0.2580696232198201
$ws.Cells.Item(2, 9).Value = 0.2737310694935027

$ws.Cells.Item(3, 1).Value = "model_2_8_1"
$ws.Cells.Item(3, 2).Value = 0.9820914312903619
$ws.Cells.Item(3, 3).Value = 0.9984512621611872
$ws.Cells.Item(3, 4).Value = 0.9926154884046207
$ws.Cells.Item(3, 5).Value = 0.9966685168192737
$ws.Cells.Item(3, 6).Value = 2.159778718311473
$ws.Cells.Item(3, 7).Value = 0.2361137515274295
$ws.Cells.Item(3, 8).Value = 0.5684926198111886
$ws.Cells.Item(3, 9).Value = 0.3925277627588912

$ws.Cells.Item(4, 1).Value = "model_2_8_2"
$ws.Cells.Item(4, 2).Value = 0.9848122297208539
$ws.Cells.Item(4, 3).Value = 0.9984335256785457
$ws.Cells.Item(4, 4).Value = 0.9871721694005128
$ws.Cells.Item(4, 5).Value = 0.9949826800241519
$ws.Cells.Item(4, 6).Value = 1.831649617529162
$ws.Cells.Item(4, 7).Value = 0.2388177775739376
$ws.Cells.Item(4, 8).Value = 0.9875435808862066
$ws.Cells.Item(4, 9).Value = 0.5911593360455762

$ws.Cells.Item(5, 1).Value = "model_2_8_24"
$ws.Cells.Item(5, 2).Value = 0.9849063877888722
$ws.Cells.Item(5, 3).Value = 0.9849033250704754
$ws.Cells.Item(5, 4).Value = 0.869570277643789
$ws.Cells.Item(5, 5).Value = 0.9495544357999471
$ws.Cells.Item(5, 6).Value = 1.820294126492413
$ws.Cells.Item(5, 7).Value = 2.301572586314831
$ws.Cells.Item(5, 8).Value = 10.04106143051149
$ws.Cells.Item(5, 9).Value = 5.943684353897952

$ws.Cells.Item(6, 1).Value = "model_2_8_23"
$ws.Cells.Item(6, 2).Value = 0.9851506795086122
$ws.Cells.Item(6, 3).Value = 0.9852343102845326
$ws.Cells.Item(6, 4).Value = 0.8718521992748797
$ws.Cells.Item(6, 5).Value = 0.9504828234309309
$ws.Cells.Item(6, 6).Value = 1.790832472358644
$ws.Cells.Item(6, 7).Value = 2.251112037968547
$ws.Cells.Item(6, 8).Value = 9.865388931456271
$ws.Cells.Item(6, 9).Value = 5.834298263681021

$ws.Cells.Item(7, 1).Value = "model_2_8_22"
$ws.Cells.Item(7, 2).Value = 0.9854148023075233
$ws.Cells.Item(7, 3).Value = 0.9855966588998485
$ws.Cells.Item(7, 4).Value = 0.8743617407183384
$ws.Cells.Item(7, 5).Value = 0.9515026582332878
$ws.Cells.Item(7, 6).Value = 1.758979184172533
$ws.Cells.Item(7, 7).Value = 2.195869963565173
$ws.Cells.Item(7, 8).Value = 9.672193244606882
$ws.Cells.Item(7, 9).Value = 5.714137526965085

$ws.Cells.Item(8, 1).Value = "model_2_8_21"
$ws.Cells.Item(8, 2).Value = 0.9856993042943372
$ws.Cells.Item(8, 3).Value = 0.9859925301719603
$ws.Cells.Item(8, 4).Value = 0.8771188280219026
$ws.Cells.Item(8, 5).Value = 0.952621531184993
$ws.Cells.Item(8, 6).Value = 1.724668159857822
$ws.Cells.Item(8, 7).Value = 2.135517172513132
$ws.Cells.Item(8, 8).Value = 9.45994037398615
$ws.Cells.Item(8, 9).Value = 5.582307746438165

$ws.Cells.Item(9, 1).Value = "model_2_8_20"
$ws.Cells.Item(9, 2).Value = 0.9860043360477887
$ws.Cells.Item(9, 3).Value = 0.9864242348283947
$ws.Cells.Item(9, 4).Value = 0.8801428860467967
$ws.Cells.Item(9, 5).Value = 0.9538470964536024
$ws.Cells.Item(9, 6).Value = 1.687881239574285
$ws.Cells.Item(9, 7).Value = 2.069701381468268
$ws.Cells.Item(9, 8).Value = 9.227134907188749
$ws.Cells.Item(9, 9).Value = 5.437907079556405

$ws.Cells.Item(10, 1).Value = "model_2_8_19"
$ws.Cells.Item(10, 2).Value = 0.986329581965286
$ws.Cells.Item(10, 3).Value = 0.9868937304496058
$ws.Cells.Item(10, 4).Value = 0.88345463149779
$ws.Cells.Item(10, 5).Value = 0.9551870535539433
$ws.Cells.Item(10, 6).Value = 1.648656485088447
$ws.Cells.Item(10, 7).Value = 1.998124146334127
$ws.Cells.Item(10, 8).Value = 8.972181979934756
$ws.Cells.Item(10, 9).Value = 5.280028340791453

$ws.Cells.Item(11, 1).Value = "model_2_8_18"
$ws.Cells.Item(11, 2).Value = 0.9866740511495466
$ws.Cells.Item(11, 3).Value = 0.98740305277367
$ws.Cells.Item(11, 4).Value = 0.887074553450207
$ws.Cells.Item(11, 5).Value = 0.9566489405741599
$ws.Cells.Item(11, 6).Value = 1.607113398907581
$ws.Cells.Item(11, 7).Value = 1.920475107447311
$ws.Cells.Item(11, 8).Value = 8.6935042518736
$ws.Cells.Item(11, 9).Value = 5.107783364508309

$ws.Cells.Item(12, 1).Value = "model_2_8_3"
$ws.Cells.Item(12, 2).Value = 0.9867732685939844
$ws.Cells.Item(12, 3).Value = 0.9981450584133535
$ws.Cells.Item(12, 4).Value = 0.9807272884942881
$ws.Cells.Item(12, 5).Value = 0.992803430701472
$ws.Cells.Item(12, 6).Value = 1.595147745568317
$ws.Cells.Item(12, 7).Value = 0.2827962266506355
$ws.Cells.Item(12, 8).Value = 1.48369924174851
$ws.Cells.Item(12, 9).Value = 0.8479266119766844

$ws.Cells.Item(13, 1).Value = "model_2_8_17"
$ws.Cells.Item(13, 2).Value = 0.9870360300214651
$ws.Cells.Item(13, 3).Value = 0.9879538110796091
$ws.Cells.Item(13, 4).Value = 0.8910233508682848
$ws.Cells.Item(13, 5).Value = 0.9582403793217309
$ws.Cells.Item(13, 6).Value = 1.563458639182026
$ws.Cells.Item(13, 7).Value = 1.836508921214091
$ws.Cells.Item(13, 8).Value = 8.389508224470607
$ws.Cells.Item(13, 9).Value = 4.920274120948014

$ws.Cells.Item(14, 1).Value = "model_2_8_16"
$ws.Cells.Item(14, 2).Value = 0.987412579787314
$ws.Cells.Item(14, 3).Value = 0.9885471631385961
$ws.Cells.Item(14, 4).Value = 0.8953201832777267
$ws.Cells.Item(14, 5).Value = 0.9599680239019068
$ws.Cells.Item(14, 6).Value = 1.518046625310262
$ws.Cells.Item(14, 7).Value = 1.746049078939344
$ws.Cells.Item(14, 8).Value = 8.058718912031622
$ws.Cells.Item(14, 9).Value = 4.716716598634134

$ws.Cells.Item(15, 1).Value = "model_2_8_15"
$ws.Cells.Item(15, 2).Value = 0.9877994077921818
$ws.Cells.Item(15, 3).Value = 0.9891837751803197
$ws.Cells.Item(15, 4).Value = 0.8999830868213636
$ws.Cells.Item(15, 5).Value = 0.9618378708377556
$ws.Cells.Item(15, 6).Value = 1.471395052752667
$ws.Cells.Item(15, 7).Value = 1.648994010178247
$ws.Cells.Item(15, 8).Value = 7.69974781188361
$ws.Cells.Item(15, 9).Value = 4.496404264873436

$ws.Cells.Item(16, 1).Value = "model_2_8_4"
$ws.Cells.Item(16, 2).Value = 0.9881365983889044
$ws.Cells.Item(16, 3).Value = 0.9976550481352612
$ws.Cells.Item(16, 4).Value = 0.9736217572252758
$ws.Cells.Item(16, 5).Value = 0.9902829723762228
$ws.Cells.Item(16, 6).Value = 1.43072976639595
$ws.Cells.Item(16, 7).Value = 0.3575010360430797
$ws.Cells.Item(16, 8).Value = 2.030714712453264
$ws.Cells.Item(16, 9).Value = 1.144896404068338

$ws.Cells.Item(17, 1).Value = "model_2_8_14"
$ws.Cells.Item(17, 2).Value = 0.9881903151060153
$ws.Cells.Item(17, 3).Value = 0.9898634247606642
$ws.Cells.Item(17, 4).Value = 0.9050267253883917
$ws.Cells.Item(17, 5).Value = 0.9638542102746335
$ws.Cells.Item(17, 6).Value = 1.424251514319275
$ws.Cells.Item(17, 7).Value = 1.545377627781211
$ws.Cells.Item(17, 8).Value = 7.31146603247052
$ws.Cells.Item(17, 9).Value = 4.258831638753308

$ws.Cells.Item(18, 1).Value = "model_2_8_13"
$ws.Cells.Item(18, 2).Value = 0.9885768238697205
$ws.Cells.Item(18, 3).Value = 0.9905848553337639
$ws.Cells.Item(18, 4).Value = 0.9104623483037559
$ws.Cells.Item(18, 5).Value = 0.9660197284585795
$ws.Cells.Item(18, 6).Value = 1.377638442340931
$ws.Cells.Item(18, 7).Value = 1.435391499198141
$ws.Cells.Item(18, 8).Value = 6.893007550613069
$ws.Cells.Item(18, 9).Value = 4.003682216755402

$ws.Cells.Item(19, 1).Value = "model_2_8_12"
$ws.Cells.Item(19, 2).Value = 0.9889473421046625
$ws.Cells.Item(19, 3).Value = 0.9913453313798886
$ws.Cells.Item(19, 4).Value = 0.9162940131616436
$ws.Cells.Item(19, 5).Value = 0.9683337880138346
$ws.Cells.Item(19, 6).Value = 1.332953832892301
$ws.Cells.Item(19, 7).Value = 1.319452669722084
$ws.Cells.Item(19, 8).Value = 6.444059994623601
$ws.Cells.Item(19, 9).Value = 3.731031096866782

$ws.Cells.Item(20, 1).Value = "model_2_8_5"
$ws.Cells.Item(20, 2).Value = 0.9890335081339107
$ws.Cells.Item(20, 3).Value = 0.9970195699733885
$ws.Cells.Item(20, 4).Value = 0.9661329655034008
$ws.Cells.Item(20, 5).Value = 0.9875450392125561
$ws.Cells.Item(20, 6).Value = 1.322562184110717
$ws.Cells.Item(20, 7).Value = 0.4543832384747895
$ws.Cells.Item(20, 8).Value = 2.607235281241179
$ws.Cells.Item(20, 9).Value = 1.467489892018395

$ws.Cells.Item(21, 1).Value = "model_2_8_11"
$ws.Cells.Item(21, 2).Value = 0.9892865747048025
$ws.Cells.Item(21, 3).Value = 0.9921404287675534
$ws.Cells.Item(21, 4).Value = 0.9225190941012169
$ws.Cells.Item(21, 5).Value = 0.9707925010006683
$ws.Cells.Item(21, 6).Value = 1.292042280315483
$ws.Cells.Item(21, 7).Value = 1.198235622958957
$ws.Cells.Item(21, 8).Value = 5.964825514974458
$ws.Cells.Item(21, 9).Value = 3.441336370634467

$ws.Cells.Item(22, 1).Value = "model_2_8_6"
$ws.Cells.Item(22, 2).Value = 0.9895699951163887
$ws.Cells.Item(22, 3).Value = 0.9962840609876819
$ws.Cells.Item(22, 4).Value = 0.9584823059400455
$ws.Cells.Item(22, 5).Value = 0.9846888137836446
$ws.Cells.Item(22, 6).Value = 1.257861694295267
$ws.Cells.Item(22, 7).Value = 0.5665156998540896
$ws.Cells.Item(22, 8).Value = 3.196217158008341
$ws.Cells.Item(22, 9).Value = 1.804021015462724

$ws.Cells.Item(23, 1).Value = "model_2_8_10"
$ws.Cells.Item(23, 2).Value = 0.9895744836232155
$ws.Cells.Item(23, 3).Value = 0.9929633593387591
$ws.Cells.Item(23, 4).Value = 0.9291236967414739
$ws.Cells.Item(23, 5).Value = 0.9733869681524723
$ws.Cells.Item(23, 6).Value = 1.257320379035583
$ws.Cells.Item(23, 7).Value = 1.07277525158782
$ws.Cells.Item(23, 8).Value = 5.456373763050747
$ws.Cells.Item(23, 9).Value = 3.135646582812377

$ws.Cells.Item(24, 1).Value = "model_2_8_9"
$ws.Cells.Item(24, 2).Value = 0.9897850421973394
$ws.Cells.Item(24, 3).Value = 0.99380449419859
$ws.Cells.Item(24, 4).Value = 0.9360800608667462
$ws.Cells.Item(24, 5).Value = 0.9761021150600636
$ws.Cells.Item(24, 6).Value = 1.231926952306512
$ws.Cells.Item(24, 7).Value = 0.9445395345297403
$ws.Cells.Item(24, 8).Value = 4.920841843998556
$ws.Cells.Item(24, 9).Value = 2.815737856463588

$ws.Cells.Item(25, 1).Value = "model_2_8_7"
$ws.Cells.Item(25, 2).Value = 0.9898310852833329
$ws.Cells.Item(25, 3).Value = 0.9954848334146408
$ws.Cells.Item(25, 4).Value = 0.950841747230103
$ws.Cells.Item(25, 5).Value = 0.9817920262340463
$ws.Cells.Item(25, 6).Value = 1.226374142427268
$ws.Cells.Item(25, 7).Value = 0.6883624164937069
$ws.Cells.Item(25, 8).Value = 3.784421426054221
$ws.Cells.Item(25, 9).Value = 2.145331319116657

$ws.Cells.Item(26, 1).Value = "model_2_8_8"
$ws.Cells.Item(26, 2).Value = 0.9898847843448165
$ws.Cells.Item(26, 3).Value = 0.9946507971294917
$ws.Cells.Item(26, 4).Value = 0.9433423501551944
$ws.Cells.Item(26, 5).Value = 0.9789148069245747
$ws.Cells.Item(26, 6).Value = 1.219898019624474
$ws.Cells.Item(26, 7).Value = 0.8155159161121485
$ws.Cells.Item(26, 8).Value = 4.361758442193086
$ws.Cells.Item(26, 9).Value = 2.484336019799967
